$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.036.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.320.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +25.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.60%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.67%  "
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.642.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.867"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.320.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.921.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  +4.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "43.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.21%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "178.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0939"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.76%  "
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +23.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0358"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.246"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.26%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  +5.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("E50").Value = "  +4.58%  "
$ws.Range("E51").Value = "  +11.24%  "
